# Auto-generated Excel COM-interop script applying the Ratios_add.xlsx edit:
# - Lab # column (A) switched from bare numeric IDs to disambiguated string labels
#   (shared strings), since several rows share the same lab number.
# - Age-calculation columns (N:S) recomputed with updated (numba-based) results.
# - Column S width trimmed by one character unit to match column Q.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Lab # (column A): replace bare numeric IDs with string labels ---
$ws.Range("A2").Value = "10815a"
$ws.Range("A3").Value = "11069b"
$ws.Range("A4").Value = "10815a"
$ws.Range("A5").Value = "11070a"
$ws.Range("A6").Value = "10815a"
$ws.Range("A7").Value = "11071a"
$ws.Range("A8").Value = "10815a"
$ws.Range("A9").Value = "11072a"
$ws.Range("A10").Value = "10815a"
$ws.Range("A11").Value = "11074b"
$ws.Range("A12").Value = "10815a"

# --- Recalculated ratio/age columns (N:S) for each row ---
# Row 2
$ws.Range("N2").Value = 1.318826771892574
$ws.Range("O2").Value = 0.5335004453826736
$ws.Range("P2").Value = 0.1702371670050599
$ws.Range("Q2").Value = 1.612895259779135
$ws.Range("R2").Value = 0.2257985771763353
$ws.Range("S2").Value = 1.726512574742643
# Row 3
$ws.Range("N3").Value = 0.03525043890623247
$ws.Range("O3").Value = 5.705288904693922
$ws.Range("P3").Value = 0.3294613250209388
$ws.Range("Q3").Value = 1.213265190703177
$ws.Range("R3").Value = 0.01157264535611715
$ws.Range("S3").Value = 2.453371522281945
# Row 4
$ws.Range("N4").Value = 1.318759177757272
$ws.Range("O4").Value = 0.6996313538309004
$ws.Range("P4").Value = 0.1720259765033036
$ws.Range("Q4").Value = 1.288866266384345
$ws.Range("R4").Value = 0.2271470050020687
$ws.Range("S4").Value = 1.004207089032775
# Row 5
$ws.Range("N5").Value = 0.005921825770778144
$ws.Range("O5").Value = 2.813032265449945
$ws.Range("P5").Value = 0.2259939097674143
$ws.Range("Q5").Value = 1.381875981856847
$ws.Range("R5").Value = 0.001274691673562648
$ws.Range("S5").Value = 2.97747271498213
# Row 6
$ws.Range("N6").Value = 1.313374031600917
$ws.Range("O6").Value = 0.5419107723895897
$ws.Range("P6").Value = 0.1815114725115821
$ws.Range("Q6").Value = 1.276965247862709
$ws.Range("R6").Value = 0.2393293430170869
$ws.Range("S6").Value = 1.28265146190901
# Row 7
$ws.Range("N7").Value = 0.01783305396177997
$ws.Range("O7").Value = 1.001748480524192
$ws.Range("P7").Value = 0.249544100283467
$ws.Range("Q7").Value = 0.5610059207786317
$ws.Range("R7").Value = 0.004458367135931112
$ws.Range("S7").Value = 0.6951876357866389
# Row 8
$ws.Range("N8").Value = 1.313892579368417
$ws.Range("O8").Value = 0.4004394967171863
$ws.Range("P8").Value = 0.1792768137488384
$ws.Range("Q8").Value = 1.36931034548719
$ws.Range("R8").Value = 0.2374191339979786
$ws.Range("S8").Value = 1.300845587375957
# Row 9
$ws.Range("N9").Value = 0.1438931916272418
$ws.Range("O9").Value = 0.365815313729743
$ws.Range("P9").Value = 0.01854280137926809
$ws.Range("Q9").Value = 0.2678043872163008
$ws.Range("R9").Value = 0.002664781889775179
$ws.Range("S9").Value = 0.3003130595442578
# Row 10
$ws.Range("N10").Value = 1.312935307459668
$ws.Range("O10").Value = 0.5991393735161512
$ws.Range("P10").Value = 0.1781077787529703
$ws.Range("Q10").Value = 1.347426895188685
$ws.Range("R10").Value = 0.2340798179609642
$ws.Range("S10").Value = 1.249551913334694
# Row 11
$ws.Range("N11").Value = 0.1434984853044898
$ws.Range("O11").Value = 0.3364171721037157
$ws.Range("P11").Value = 0.0009631928586580701
$ws.Range("Q11").Value = 0.255338915736973
$ws.Range("R11").Value = 0.0001387130736508046
$ws.Range("S11").Value = 0.2564753741715639
# Row 12
$ws.Range("N12").Value = 1.318059551970325
$ws.Range("O12").Value = 0.4537665522382115
$ws.Range("P12").Value = 0.1784094620679704
$ws.Range("Q12").Value = 1.295662586723902
$ws.Range("R12").Value = 0.2343682381076229
$ws.Range("S12").Value = 1.283193429135786

# --- Column S width: narrow by one character unit (20.71 -> 19.71) ---
$ws.Columns.Item(19).ColumnWidth = $ws.Columns.Item(19).ColumnWidth - 1

